$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 = "time_taken", matching the style of the other header cells (e.g. E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# New data column F2:F5 with time_taken values (plain text, same as other data cells)
$ws.Range("F2").Value = "2021-10-05 13:42:32.881278"
$ws.Range("F3").Value = "2021-10-05 13:42:32.881321"
$ws.Range("F4").Value = "2021-10-05 13:42:32.881326"
$ws.Range("F5").Value = "2021-10-05 13:42:32.881330"
